$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.212486028671265
$ws.Range("B1").Value = 2.673704385757446
$ws.Range("C1").Value = 2.812166929244995
$ws.Range("D1").Value = 2.554224729537964
$ws.Range("E1").Value = 0.8132590055465698
